$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at EF, shifting old EF (nom) -> EG and old EG (url) -> EH
$ws.Columns("EF").Insert()

# Header row: new timestamp column
$ws.Range("EF1").Value = "2026-02-03 05:42:58"

# Rows 2-80 have numeric price history in EE; copy EE's value into the new EF cell
for ($r = 2; $r -le 80; $r++) {
    $ee = $ws.Cells.Item($r, 135)
    $ef = $ws.Cells.Item($r, 136)
    $ef.Value = $ee.Value2
}

# Special-case text corrections (observed scrape anomalies) on the shifted "nom" column (EG)
$ws.Range("EG19").Value = "Apple iPhone 16 128 Go Rose809€00"
$ws.Range("EG76").Value = "Apple iPhone 17 Pro Max 1 To Bleu Intense"
